$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for all data rows (2 through 260) is being updated
# from the old "changed" date (45172 = 2023-09-03) to the new one
# (45175 = 2023-09-06).
$ws.Range("C2:C260").Value = 45175
